# Added column mapping listing api
# Clear the "Source" column (A) values for the rows that are list/array
# type hierarchy levels which only have a mapping in the "Target" column (B).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("A6").ClearContents()
$ws.Range("A10").ClearContents()
$ws.Range("A12:A18").ClearContents()
$ws.Range("A20").ClearContents()
